$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold text that looks numeric (e.g. "8.10", "517.48").
# Force the Text number format first so Excel does not silently coerce these
# assignments into Double values (which would drop the trailing zeros / dots).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.528.02"
$ws.Range("E2").Value = "  +2.43%  "

$ws.Range("D3").Value = "3.059.09"
$ws.Range("E3").Value = "  +2.67%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "517.48"
$ws.Range("E5").Value = "  +2.73%  "

$ws.Range("D6").Value = "141.58"
$ws.Range("E6").Value = "  +3.24%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "0.433"
$ws.Range("E8").Value = "  +1.63%  "

$ws.Range("D9").Value = "7.25"
$ws.Range("E9").Value = "  +1.38%  "

$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("D11").Value = "0.374"
$ws.Range("E11").Value = "  +2.97%  "

$ws.Range("D12").Value = "3.592.37"
$ws.Range("E12").Value = "  +3.03%  "

$ws.Range("E13").Value = "  +3.19%  "

$ws.Range("D14").Value = "25.54"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").Value = "0.0000162"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").Value = "57.607.25"
$ws.Range("E16").Value = "  +2.61%  "

$ws.Range("D17").Value = "3.063.35"
$ws.Range("E17").Value = "  +2.78%  "

$ws.Range("D18").Value = "6.06"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -0.37%  "

$ws.Range("D20").Value = "8.10"
$ws.Range("E20").Value = "  +1.89%  "

$ws.Range("D21").Value = "330.09"
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("D23").Value = "0.496"
$ws.Range("E23").Value = "  +0.73%  "

$ws.Range("D24").Value = "65.77"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("E25").Value = "  +3.70%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "0.0₃0896"
$ws.Range("E27").Value = "  -2.49%  "

$ws.Range("D28").Value = "6.30"
$ws.Range("E28").Value = "  -0.87%  "

$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  +1.96%  "

$ws.Range("E30").Value = "  +2.83%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "20.66"
$ws.Range("E31").Value = "  +2.74%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("D33").Value = "154.71"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("D34").Value = "27.36"
$ws.Range("E34").Value = "  +6.20%  "

$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  +3.06%  "

$ws.Range("D37").Value = "1.27"
$ws.Range("E37").Value = "  +2.96%  "

$ws.Range("D38").Value = "0.0670"
$ws.Range("E38").Value = "  +1.72%  "

$ws.Range("D39").Value = "3.103.60"
$ws.Range("E39").Value = "  +2.95%  "

$ws.Range("D40").Value = "3.89"
$ws.Range("E40").Value = "  +3.09%  "

$ws.Range("D41").Value = "36.59"
$ws.Range("E41").Value = "  -0.71%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "0.655"
$ws.Range("E43").Value = "  +0.84%  "

$ws.Range("D44").Value = "2.255.13"
$ws.Range("E44").Value = "  +4.17%  "

$ws.Range("D45").Value = "0.0257"
$ws.Range("E45").Value = "  +8.93%  "

$ws.Range("D46").Value = "20.68"
$ws.Range("E46").Value = "  +6.44%  "

$ws.Range("E47").Value = "  +0.71%  "

$ws.Range("D48").Value = "5.86"
$ws.Range("E48").Value = "  +0.90%  "

$ws.Range("D49").Value = "0.919"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").Value = "261.10"
$ws.Range("E50").Value = "  +15.62%  "

$ws.Range("D51").Value = "0.713"
$ws.Range("E51").Value = "  +6.38%  "
